# The document has a BTEC logo in the primary header and a Pearson Edexcel
# logo in each of the two footers. The source assets backing those three
# inline pictures were swapped (the BTEC logo now comes from "image1.jpg"
# instead of "image2.jpg", and the two Pearson logos now come from
# "image2.png" instead of "image1.png"), so the shape names need to be
# updated to match.

$d = $word.ActiveDocument
$sec = $d.Sections.First

# --- Header: BTec_Logo-Orange -> rename image2.jpg to image1.jpg ---
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $header = $sec.Headers.Item($i)
    if ($header.Exists) {
        for ($j = 1; $j -le $header.Range.InlineShapes.Count; $j++) {
            $shape = $header.Range.InlineShapes.Item($j)
            if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
                $shape.Name = "image1.jpg"
            }
        }
    }
}

# --- Footers: Pearson Edexcel logo -> rename image1.png to image2.png ---
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $footer = $sec.Footers.Item($i)
    if ($footer.Exists) {
        for ($j = 1; $j -le $footer.Range.InlineShapes.Count; $j++) {
            $shape = $footer.Range.InlineShapes.Item($j)
            if ($shape.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shape.Name = "image2.png"
            }
        }
    }
}
